$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) labels to reflect revised quantile naming ---
$ws.Range("A1").Value = "beta-mean"
$ws.Range("B1").Value = "beta-Q15.5"
$ws.Range("C1").Value = "beta-Q83.5"
$ws.Range("D1").Value = "beta-Q2.5"
$ws.Range("E1").Value = "beta-Q97.5"
$ws.Range("F1").Value = "gamma-inv-mean"
$ws.Range("G1").Value = "gamma-inv-Q15.5"
$ws.Range("H1").Value = "gamma-inv-Q83.5"
$ws.Range("I1").Value = "gamma-inv-Q2.5"
$ws.Range("J1").Value = "gamma-inv-Q97.5"
$ws.Range("K1").Value = "R_0-mean"
$ws.Range("L1").Value = "R_0-Q15.5"
$ws.Range("M1").Value = "R_0-Q83.5"
$ws.Range("N1").Value = "R_0-Q2.5"
$ws.Range("O1").Value = "R_0-Q97.5"
$ws.Range("P1").Value = "t_c-mean"
$ws.Range("Q1").Value = "t_c-Q15.5"
$ws.Range("R1").Value = "t_c-Q83.5"
$ws.Range("S1").Value = "t_c-Q2.5"
$ws.Range("T1").Value = "t_c-Q97.5"
$ws.Range("U1").Value = "I_peak-mean"
$ws.Range("V1").Value = "I_peak-Q15.5"
$ws.Range("W1").Value = "I_peak-Q83.5"
$ws.Range("X1").Value = "I_peak-Q2.5"
$ws.Range("Y1").Value = "I_peak-Q97.5"
$ws.Range("Z1").Value = "T_end-mean"
$ws.Range("AA1").Value = "T_end-Q15.5"
$ws.Range("AB1").Value = "T_end-Q83.5"
$ws.Range("AC1").Value = "T_end-Q2.5"
$ws.Range("AD1").Value = "T_end-Q97.5"

# --- Update simulation result data rows (row 2-4) with recomputed values ---

# Row 2
$ws.Range("A2").Value = 0.3300499060118708
$ws.Range("B2").Value = 0.2953613910451203
$ws.Range("C2").Value = 0.3633580143913455
$ws.Range("D2").Value = 0.271946956947787
$ws.Range("E2").Value = 0.4069920332739979
$ws.Range("K2").Value = 2.310349342083096
$ws.Range("L2").Value = 2.067529737315842
$ws.Range("M2").Value = 2.543506100739418
$ws.Range("N2").Value = 1.903628698634509
$ws.Range("O2").Value = 2.848944232917985
$ws.Range("P2").Value = 64.48788999999999
$ws.Range("Q2").Value = 46.20353251271484
$ws.Range("R2").Value = 87.6772320006519
$ws.Range("S2").Value = 53.56114539921523
$ws.Range("T2").Value = 74.57462299720029
$ws.Range("U2").Value = 0.2030477690805852
$ws.Range("V2").Value = 0.1643481661970219
$ws.Range("W2").Value = 0.2419927606527416
$ws.Range("X2").Value = 0.1356656996107907
$ws.Range("Y2").Value = 0.2800230314464088
$ws.Range("Z2").Value = 0.8571281783494924
$ws.Range("AA2").Value = 0.81085013287319
$ws.Range("AB2").Value = 0.8980873652968309
$ws.Range("AC2").Value = 0.7700111195685573
$ws.Range("AD2").Value = 0.9302255099274735

# Row 3
$ws.Range("F3").Value = 6.997911609327875
$ws.Range("G3").Value = 6.210342225982567
$ws.Range("H3").Value = 7.747438010729122
$ws.Range("I3").Value = 5.693943946463319
$ws.Range("J3").Value = 8.779566248620817
$ws.Range("K3").Value = 2.309310831078199
$ws.Range("L3").Value = 2.049412934574248
$ws.Range("M3").Value = 2.55665454354061
$ws.Range("N3").Value = 1.879001502332895
$ws.Range("O3").Value = 2.89725686204487
$ws.Range("P3").Value = 63.53957
$ws.Range("Q3").Value = 56.39647523232526
$ws.Range("R3").Value = 72.30353979352127
$ws.Range("S3").Value = 59.40480311769122
$ws.Range("T3").Value = 67.82739293048661
$ws.Range("U3").Value = 0.2026521904241902
$ws.Range("V3").Value = 0.1614428624796646
$ws.Range("W3").Value = 0.2444269455347663
$ws.Range("X3").Value = 0.1310788822784472
$ws.Range("Y3").Value = 0.2858611211784558
$ws.Range("Z3").Value = 0.8558996011535711
$ws.Range("AA3").Value = 0.8062071331525196
$ws.Range("AB3").Value = 0.8998128159466621
$ws.Range("AC3").Value = 0.7623227833969787
$ws.Range("AD3").Value = 0.9341038326876447

# Row 4
$ws.Range("A4").Value = 0.3300191363416627
$ws.Range("B4").Value = 0.2953936898295382
$ws.Range("C4").Value = 0.3631343751871358
$ws.Range("D4").Value = 0.2719975722997606
$ws.Range("E4").Value = 0.406778017689316
$ws.Range("F4").Value = 7.00047510005597
$ws.Range("G4").Value = 6.213335460656964
$ws.Range("H4").Value = 7.754827506595099
$ws.Range("I4").Value = 5.69610889909607
$ws.Range("J4").Value = 8.781418602459933
$ws.Range("K4").Value = 2.310420417988074
$ws.Range("L4").Value = 1.957020368334278
$ws.Range("M4").Value = 2.649412775043653
$ws.Range("N4").Value = 1.720654360957462
$ws.Range("O4").Value = 3.113987979913325
$ws.Range("P4").Value = 65.38379999999999
$ws.Range("Q4").Value = 45.5884797976407
$ws.Range("R4").Value = 94.57954781029983
$ws.Range("S4").Value = 53.75743574891678
$ws.Range("T4").Value = 76.64907939477872
$ws.Range("U4").Value = 0.2013060712975479
$ws.Range("V4").Value = 0.1448516681886182
$ws.Range("W4").Value = 0.2578206928040442
$ws.Range("X4").Value = 0.102701757916482
$ws.Range("Y4").Value = 0.3125445047127944
$ws.Range("Z4").Value = 0.8487886733785018
$ws.Range("AA4").Value = 0.776808401230911
$ws.Range("AB4").Value = 0.9082350849830558
$ws.Range("AC4").Value = 0.7052318602969566
$ws.Range("AD4").Value = 0.9518547260593895
